$d = $word.ActiveDocument

# Unicode curly quote characters used throughout the document.
$lq = [char]0x201C   # “
$rq = [char]0x201D   # ”

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# 1) Rewrite the "Into the project, import dirs ..." bullet: it now only
#    mentions "data" and "src" (not "bin"/"lib"), and gains a brand-new
#    bullet right after it ("In the project, create a dir called "lib".").
# ---------------------------------------------------------------------------
$find = $d.Content
$find.Find.Execute("Into the project, import") | Out-Null
$targetPara = $find.Paragraphs(1).Range
$targetPara.MoveEnd(1, -1) | Out-Null

$body =
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Into the project, import </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>dirs</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    ('<w:r><w:t xml:space="preserve"> {0}data{1}</w:t></w:r>' -f $lq, $rq) +
    '<w:r><w:t xml:space="preserve"> and</w:t></w:r>' +
    ('<w:r><w:t xml:space="preserve"> {0}</w:t></w:r>' -f $lq) +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>src</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    ('<w:r><w:t xml:space="preserve">{0}. You can easily do this by selecting </w:t></w:r>' -f $rq) +
    '<w:r><w:t>the</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>dirs</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> and dragging them into the icon for the Java workspace.</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">In the project, create a </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>dir</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    ('<w:r><w:t xml:space="preserve"> called {0}lib{1}.</w:t></w:r>' -f $lq, $rq) +
  '</w:p>'

$targetPara.InsertXML($pkgHeader + $body + $pkgFooter)

# ---------------------------------------------------------------------------
# 2) Move the stale "_GoBack" bookmark so it again sits at the very last
#    edit location - right after the text that was just typed ("lib".).
# ---------------------------------------------------------------------------
$goBackSpot = $d.Content
$goBackSpot.Find.Execute("called " + $lq + "lib" + $rq + ".") | Out-Null
$goBackSpot.Collapse(0) | Out-Null
$d.Bookmarks.Add("_GoBack", $goBackSpot) | Out-Null

# ---------------------------------------------------------------------------
# 3) Drop the stale <w:lastRenderedPageBreak/> cached before the "Execu"
#    run of the "Execution:" heading - it no longer corresponds to a real
#    page break once the document above it changed length.
# ---------------------------------------------------------------------------
$execRangeFinder = $d.Content
$execRangeFinder.Find.Execute("Execution:") | Out-Null
$execPara = $execRangeFinder.Paragraphs(1).Range
$execPara.MoveEnd(1, -1) | Out-Null

$execBody =
  '<w:p><w:pPr><w:rPr><w:sz w:val="36"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:sz w:val="36"/></w:rPr><w:t>Execu</w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="36"/></w:rPr><w:t>tion:</w:t></w:r>' +
  '</w:p>'

$execPara.InsertXML($pkgHeader + $execBody + $pkgFooter)
